$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns to hold the new inventory-summary fields.
# (Each Insert() shifts everything at/after that column one to the right.)
$ws.Columns("U").Insert()   # new: Total Stall Count
$ws.Columns("Y").Insert()   # new: Lift Per Tower
$ws.Columns("AA").Insert()  # new: Flier Frequency

# Header row updates
$ws.Range("U1").Value2 = "Total Stall Count"
$ws.Range("Y1").Value2 = "Lift Per Tower"
$ws.Range("Z1").Value2 = "Flier Allowed(Y/N)"
$ws.Range("AA1").Value2 = "Flier Frequency"

# Row 2 (Jelly Beans) new data
$ws.Range("E2").Value2 = "ABC"
$ws.Range("U2").Value2 = 1
$ws.Range("Y2").Value2 = 2
$ws.Range("AA2").Value2 = 2

# Row 3 (Choco Pie) new data
$ws.Range("E3").Value2 = "ADF"
$ws.Range("U3").Value2 = 2
$ws.Range("Y3").Value2 = 3
$ws.Range("AA3").Value2 = 3

# Column width tweaks to match the post-edit layout
$ws.Columns("Z").ColumnWidth = 9.25
$ws.Columns("AA").ColumnWidth = 12.084

# View / layout tweaks
$excel.ActiveWindow.Zoom = 75
$ws.Rows(1).RowHeight = 55.75
$ws.Range("E5").Select()
